$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column H (shifts old H:M to K:P)
$ws.Range("H1:J1").EntireColumn.Insert()

# New header row values
$ws.Range("H1").Value = "excel"
$ws.Range("I1").Value = "csv"
$ws.Range("J1").Value = "pdf"

# New data values for row 2
$ws.Range("H2").Value = "SI"
$ws.Range("I2").Value = "SI"
$ws.Range("J2").Value = "NO"

# New data values for row 3
$ws.Range("H3").Value = "NO"
$ws.Range("I3").Value = "NO"
$ws.Range("J3").Value = "NO"

Write-Output "done"
